# Append 18 new transaction rows (22-39) to the "Transactions" sheet.
#
# Column B holds transaction reference numbers that are numeric-looking but
# must be stored as TEXT (matching the existing rows 1-21, which are all
# text/shared-strings). Assigning a plain numeric-looking string via
# Range.Value gets auto-converted to a Number by Excel, and prefixing with
# a literal apostrophe marks the cell with a "quotePrefix" style that the
# original rows don't have. To get clean text cells (same shape as the
# existing data) we instead enter a text formula (="1234...") for each row
# and then convert it to a static value via Copy + PasteSpecial (values
# only), which leaves a plain text cell behind with no extra formula or
# style baggage.
#
# Column C holds a status note - "YOUR PAYMENT WAS DECLINED" - for the two
# declined transactions (rows 36 and 39). That text isn't numeric-looking,
# so it can be entered directly with Value without being misread as a
# number.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transactions")

function Set-TextNumber($cellRef, $numericText) {
    $cell = $ws.Range($cellRef)
    $cell.Formula = '="' + $numericText + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

Set-TextNumber "B22" "3320036608"
Set-TextNumber "B23" "3318789888"
Set-TextNumber "B24" "3335400448"
Set-TextNumber "B25" "3386408704"
Set-TextNumber "B26" "3373721088"
Set-TextNumber "B27" "3301729792"
Set-TextNumber "B28" "3309341952"
Set-TextNumber "B29" "3353609728"
Set-TextNumber "B30" "3311197184"
Set-TextNumber "B31" "3392648960"
Set-TextNumber "B32" "3385755136"
Set-TextNumber "B33" "3340498176"
Set-TextNumber "B34" "3389296640"
Set-TextNumber "B35" "3319841000"
Set-TextNumber "B36" "3335292672"
$ws.Range("C36").Value = "YOUR PAYMENT WAS DECLINED"
Set-TextNumber "B37" "3399858176"
Set-TextNumber "B38" "3388921856"
Set-TextNumber "B39" "3397952512"
$ws.Range("C39").Value = "YOUR PAYMENT WAS DECLINED"
